$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B, D, E, F, G, H, I, J, O (rows 2-25),
# as per the updated 380 kV line power-flow results (pl_mw.xlsx).
# Each row entry: row number, then values for B,D,E,F,G,H,I,J,O in that order.
$data = @(
    @(2, 1.219803940679071, 0.246097157014546, 0.1773659978598765, 1.012249720703693, 0.5272962014676921, 0.5812333897015094, 3.961193417866525, 0.1772330055721199, 2.183072434228677),
    @(3, 1.069997499266549, 0.2444676174134202, 0.1742912681868667, 0.9831260700589439, 0.4999711970626919, 0.5723711798586777, 3.481483530427994, 0.1721769985554289, 2.106238768013412),
    @(4, 0.9775436500530645, 0.2435737801817908, 0.1724976178116897, 0.9658221175057946, 0.4835087855673947, 0.5672287518405454, 3.186307868167489, 0.1691764068403643, 2.060343181411326),
    @(5, 0.9397517899408854, 0.2432364075533258, 0.1717903946035122, 0.9589159623540269, 0.4768793569706844, 0.5652084079914204, 3.065876933603903, 0.167979719279117, 2.041962182227564),
    @(6, 0.9334695415156489, 0.24318201137892, 0.171674392540087, 0.957777978631924, 0.4757833212346156, 0.5648774770931766, 3.04587117873993, 0.1677825844652077, 2.038929463360034),
    @(7, 0.9770344427063833, 0.2435691213927882, 0.1724879839689208, 0.9657283902091791, 0.4834190585289946, 0.5672012000805466, 3.184684259858585, 0.1691601623408303, 2.060093986032228),
    @(8, 1.168250105383493, 0.2455131683178706, 0.1762862718187961, 1.002087880873745, 0.517809018632434, 0.5781156669493015, 3.795928923952232, 0.1754681284190411, 2.156314253543428),
    @(9, 1.539389405566169, 0.2501707657315677, 0.1844829058625663, 1.077981530534032, 0.5877591867547096, 0.6018910408742215, 4.989014804984947, 0.1886638921441488, 2.355184202377018),
    @(10, 1.809644099351658, 0.2541069149559831, 0.1909625588303214, 1.136556074637113, 0.6407026906574345, 0.6208072618303788, 5.861568061238756, 0.1988664766722508, 2.507549276033956),
    @(11, 1.932050532542917, 0.2560090365128502, 0.194010052836056, 1.1638183204999, 0.6651295264913699, 0.6297279495998964, 6.25753114434093, 0.2036191076019378, 2.578233826267081),
    @(12, 1.978324344473322, 0.256745331136699, 0.195178431529726, 1.174230611703237, 0.6744288402402958, 0.6331513654079117, 6.407322560669854, 0.2054348799560159, 2.605198180176444),
    @(13, 1.968361989860625, 0.2565860458044114, 0.1949261616900486, 1.17198419430666, 0.6724238672803722, 0.6324120554357648, 6.375069214559744, 0.2050431063757969, 2.599382130016579),
    @(14, 1.935859098657602, 0.2560692914184273, 0.1941058882925901, 1.164673168681375, 0.6658935957561596, 0.6300086875166357, 6.269857673435979, 0.2037681703561276, 2.580448237419887),
    @(15, 1.915939811075305, 0.2557548473068891, 0.1936053170411043, 1.160206502899541, 0.6619000535299335, 0.6285424593906441, 6.205392518134659, 0.202989326893956, 2.568876427918497),
    @(16, 1.801633628443199, 0.253984848216092, 0.19076540734887, 1.13478682902678, 0.6391132448457597, 0.6202306213541249, 5.835670441680065, 0.1985581258952323, 2.502957520220775),
    @(17, 1.731372246120657, 0.2529275540791929, 0.1890487917632697, 1.119350579198027, 0.6252221158740383, 0.6152123840409445, 5.608601612490133, 0.1958682884040428, 2.462870201641749),
    @(18, 1.690909632147168, 0.2523299281006217, 0.1880708399873754, 1.110530072250796, 0.6172645418874652, 0.6123557350702242, 5.477907426105048, 0.1943316508989312, 2.439942247904867),
    @(19, 1.677201143001923, 0.2521293870373, 0.1877413368566749, 1.107553563374239, 0.6145757763602546, 0.611393626442009, 5.433641513037287, 0.1938131718418248, 2.432201422091339),
    @(20, 1.738856893516015, 0.2530390182404858, 0.1892305555620553, 1.120987789708153, 0.6266975129549053, 0.6157435095531412, 5.632782903638144, 0.1961535407198483, 2.467124190622599),
    @(21, 1.945408141631617, 0.2562206407067436, 0.1943464327933739, 1.16681818736771, 0.6678103529508235, 0.6307133845868975, 6.30076505272973, 0.2041422140131033, 2.586004214704474),
    @(22, 2.079940855282643, 0.2583932720504407, 0.1977736478160992, 1.19728785651229, 0.6949679831682261, 0.6407613510679369, 6.736443306619492, 0.2094568644436805, 2.664851336729555),
    @(23, 2.008181005422728, 0.2572251772786558, 0.1959368213281465, 1.180978317018202, 0.6804470454097498, 0.6353743924277921, 6.503998686294949, 0.2066117615138694, 2.622663646758269),
    @(24, 1.735473297350438, 0.2529885934366121, 0.1891483522606237, 1.120247438639595, 0.6260303967946754, 0.615503299234291, 5.621851005971507, 0.1960245476697935, 2.465200592235419),
    @(25, 1.439406031741498, 0.2488203805639699, 0.1821852480318213, 1.05695698794338, 0.5685648655342277, 0.5952049802156978, 4.666914883151151, 0.1850052413102077, 2.300289579330695)
)

$cols = @(2, 4, 5, 6, 7, 8, 9, 10, 15)  # B, D, E, F, G, H, I, J, O

foreach ($entry in $data) {
    $row = $entry[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($row, $cols[$i]).Value = $entry[$i + 1]
    }
}

